# Update the "Change Start page" test step (cell G2 on the TestCases sheet)
# to use the new Configxml based approach instead of the old SetStartPage call.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$newText = "wait(3);`n" + `
    "PullConfigxml;`n" + `
    "ChangeConfigxml(Configuration/Applications/Application/General,StartPage,<StartPage value=`"http://127.0.0.1:8082/app/`" name=`"Menu`"/>);`n" + `
    "ChangeConfigxml(Configuration,WebServer,<WebServer>);`n" + `
    "ChangeConfigxml(Configuration/WebServer,Enabled,<Enabled VALUE=`"1`"/>);`n" + `
    "ChangeConfigxml(Configuration/WebServer,Port,<Port VALUE=`"8082`"/>);`n" + `
    "ChangeConfigxml(Configuration/WebServer,WebFolder,<WebFolder VALUE=`"\\auto\\ComplianceTest_JS\`"/>);`n" + `
    "ChangeConfigxml(Configuration/WebServer,Public,<Public VALUE=`"1`"/>);`n" + `
    "ChangeConfigxml(Configuration/Screen,FullScreen,<FullScreen value=`"0`"/>);`n" + `
    "PushConfigxml;"

$ws.Range("G2").Value = $newText

# Row grew a lot taller once the longer text wrapped, so reflect the new height.
$ws.Rows.Item(2).RowHeight = 332.25

# Selection in the saved view moves from J2:J5 to the edited cell G2.
$ws.Range("G2").Select()
